$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Narrow column D (separator-ish width), matching the added <col> for column 4.
# The runtime rounds ColumnWidth to the nearest 1/6-character increment, so the
# closest reachable value to the target stored width (0.5703125) is 0.5, which
# is produced by requesting a (clamped) negative ColumnWidth here.
$ws.Columns.Item(4).ColumnWidth = -0.4

# Header for the new "Level" column
$ws.Range("E1").Value = "Level"

# Level values for the first 10 molecules (picked by row number), assigned
# directly so the result does not depend on hashtable enumeration order.
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(24, 5).Value = 5
$ws.Cells.Item(26, 5).Value = 6
$ws.Cells.Item(30, 5).Value = 7
$ws.Cells.Item(32, 5).Value = 8
$ws.Cells.Item(33, 5).Value = 9
$ws.Cells.Item(49, 5).Value = 10

# Update the view: the saved file shows cell E1 selected, with the sheet
# scrolled back to the top (the old topLeftCell="A10"/activeCell C17
# selection is replaced).
$ws.Activate()
$ws.Range("E1").Select()
